# Insert a new weekly data row for Espinaca (Femacal de La Calera) at row 145.
# This shifts all existing rows from 145..273 down to 146..274 and
# leaves a blank row 145 to be populated with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(145).Insert()

$ws.Range("A145").Value = 3
$ws.Range("B145").Value = "Femacal de La Calera"
$ws.Range("C145").Value = "Coquimbo"
$ws.Range("D145").Value = 44589
$ws.Range("E145").Value = 5
$ws.Range("F145").Value = 100112012
$ws.Range("G145").Value = "Espinaca"
$ws.Range("H145").Value = "Sin especificar"
$ws.Range("I145").Value = "Primera"
$ws.Range("J145").Value = 170
$ws.Range("K145").Value = 3500
$ws.Range("L145").Value = 4000
$ws.Range("M145").Value = 3765
$ws.Range("N145").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O145").Value = "Provincia de Quillota"
$ws.Range("P145").Value = 1255
$ws.Range("Q145").Value = 3
$ws.Range("R145").Value = "Hortaliza"
